# Update "want to go" counts (column F) on the 展览 (sheet1), 演出 (sheet2)
# and 全部类型 (sheet4) worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 34
$wsExpo.Range("F5").Value = 2634
$wsExpo.Range("F7").Value = 361

# 演出 sheet
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 27

# 全部类型 sheet (combined view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 34
$wsAll.Range("F5").Value = 2634
$wsAll.Range("F7").Value = 361
$wsAll.Range("F8").Value = 27
